$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet by duplicating "2022-Q3" (so it
#    inherits identical column widths / header styles / cell formats),
#    positioned immediately before "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q4"

# The template ("2022-Q3") only has 14 data rows (rows 2-15); the new
# "2022-Q4" sheet needs 16 data rows (rows 2-17), so extend it by
# duplicating the formatting of the last existing row twice.
$newWs.Range("A15:H15").Copy($newWs.Range("A16:H16"))
$newWs.Range("A15:H15").Copy($newWs.Range("A17:H17"))

# ---------------------------------------------------------------------
# 2) Populate "2022-Q4" with the fund-holding data.
# ---------------------------------------------------------------------
$newWs.Cells.Item(2,1).Value = 0
$newWs.Cells.Item(2,2).NumberFormat = "@"
$newWs.Cells.Item(2,2).Value = "001822"
$newWs.Cells.Item(2,3).NumberFormat = "@"
$newWs.Cells.Item(2,3).Value = "华商智能生活灵活配置混合A"
$newWs.Cells.Item(2,4).NumberFormat = "@"
$newWs.Cells.Item(2,4).Value = "33.45"
$newWs.Cells.Item(2,5).NumberFormat = "@"
$newWs.Cells.Item(2,5).Value = "90.70"
$newWs.Cells.Item(2,6).NumberFormat = "@"
$newWs.Cells.Item(2,6).Value = "4.51"
$newWs.Cells.Item(2,7).NumberFormat = "@"
$newWs.Cells.Item(2,7).Value = "1.5086"
$newWs.Cells.Item(2,8).Value = 7
$newWs.Cells.Item(3,1).Value = 1
$newWs.Cells.Item(3,2).NumberFormat = "@"
$newWs.Cells.Item(3,2).Value = "000729"
$newWs.Cells.Item(3,3).NumberFormat = "@"
$newWs.Cells.Item(3,3).Value = "建信中小盘先锋股票A"
$newWs.Cells.Item(3,4).NumberFormat = "@"
$newWs.Cells.Item(3,4).Value = "31.93"
$newWs.Cells.Item(3,5).NumberFormat = "@"
$newWs.Cells.Item(3,5).Value = "91.13"
$newWs.Cells.Item(3,6).NumberFormat = "@"
$newWs.Cells.Item(3,6).Value = "3.70"
$newWs.Cells.Item(3,7).NumberFormat = "@"
$newWs.Cells.Item(3,7).Value = "1.1814"
$newWs.Cells.Item(3,8).Value = 6
$newWs.Cells.Item(4,1).Value = 2
$newWs.Cells.Item(4,2).NumberFormat = "@"
$newWs.Cells.Item(4,2).Value = "530005"
$newWs.Cells.Item(4,3).NumberFormat = "@"
$newWs.Cells.Item(4,3).Value = "建信优化配置混合A"
$newWs.Cells.Item(4,4).NumberFormat = "@"
$newWs.Cells.Item(4,4).Value = "24.05"
$newWs.Cells.Item(4,5).NumberFormat = "@"
$newWs.Cells.Item(4,5).Value = "88.67"
$newWs.Cells.Item(4,6).NumberFormat = "@"
$newWs.Cells.Item(4,6).Value = "3.60"
$newWs.Cells.Item(4,7).NumberFormat = "@"
$newWs.Cells.Item(4,7).Value = "0.8658"
$newWs.Cells.Item(4,8).Value = 6
$newWs.Cells.Item(5,1).Value = 3
$newWs.Cells.Item(5,2).NumberFormat = "@"
$newWs.Cells.Item(5,2).Value = "001933"
$newWs.Cells.Item(5,3).NumberFormat = "@"
$newWs.Cells.Item(5,3).Value = "华商新兴活力灵活配置混合"
$newWs.Cells.Item(5,4).NumberFormat = "@"
$newWs.Cells.Item(5,4).Value = "25.00"
$newWs.Cells.Item(5,5).NumberFormat = "@"
$newWs.Cells.Item(5,5).Value = "90.39"
$newWs.Cells.Item(5,6).NumberFormat = "@"
$newWs.Cells.Item(5,6).Value = "3.27"
$newWs.Cells.Item(5,7).NumberFormat = "@"
$newWs.Cells.Item(5,7).Value = "0.8175"
$newWs.Cells.Item(5,8).Value = 9
$newWs.Cells.Item(6,1).Value = 4
$newWs.Cells.Item(6,2).NumberFormat = "@"
$newWs.Cells.Item(6,2).Value = "013886"
$newWs.Cells.Item(6,3).NumberFormat = "@"
$newWs.Cells.Item(6,3).Value = "华商新能源汽车混合A"
$newWs.Cells.Item(6,4).NumberFormat = "@"
$newWs.Cells.Item(6,4).Value = "9.05"
$newWs.Cells.Item(6,5).NumberFormat = "@"
$newWs.Cells.Item(6,5).Value = "89.00"
$newWs.Cells.Item(6,6).NumberFormat = "@"
$newWs.Cells.Item(6,6).Value = "6.93"
$newWs.Cells.Item(6,7).NumberFormat = "@"
$newWs.Cells.Item(6,7).Value = "0.6272"
$newWs.Cells.Item(6,8).Value = 4
$newWs.Cells.Item(7,1).Value = 5
$newWs.Cells.Item(7,2).NumberFormat = "@"
$newWs.Cells.Item(7,2).Value = "000756"
$newWs.Cells.Item(7,3).NumberFormat = "@"
$newWs.Cells.Item(7,3).Value = "建信潜力新蓝筹股票A"
$newWs.Cells.Item(7,4).NumberFormat = "@"
$newWs.Cells.Item(7,4).Value = "15.30"
$newWs.Cells.Item(7,5).NumberFormat = "@"
$newWs.Cells.Item(7,5).Value = "90.58"
$newWs.Cells.Item(7,6).NumberFormat = "@"
$newWs.Cells.Item(7,6).Value = "3.75"
$newWs.Cells.Item(7,7).NumberFormat = "@"
$newWs.Cells.Item(7,7).Value = "0.5738"
$newWs.Cells.Item(7,8).Value = 6
$newWs.Cells.Item(8,1).Value = 6
$newWs.Cells.Item(8,2).NumberFormat = "@"
$newWs.Cells.Item(8,2).Value = "010550"
$newWs.Cells.Item(8,3).NumberFormat = "@"
$newWs.Cells.Item(8,3).Value = "华商双擎领航混合"
$newWs.Cells.Item(8,4).NumberFormat = "@"
$newWs.Cells.Item(8,4).Value = "12.41"
$newWs.Cells.Item(8,5).NumberFormat = "@"
$newWs.Cells.Item(8,5).Value = "90.98"
$newWs.Cells.Item(8,6).NumberFormat = "@"
$newWs.Cells.Item(8,6).Value = "4.61"
$newWs.Cells.Item(8,7).NumberFormat = "@"
$newWs.Cells.Item(8,7).Value = "0.5721"
$newWs.Cells.Item(8,8).Value = 7
$newWs.Cells.Item(9,1).Value = 7
$newWs.Cells.Item(9,2).NumberFormat = "@"
$newWs.Cells.Item(9,2).Value = "015385"
$newWs.Cells.Item(9,3).NumberFormat = "@"
$newWs.Cells.Item(9,3).Value = "华商智能生活灵活配置混合C"
$newWs.Cells.Item(9,4).NumberFormat = "@"
$newWs.Cells.Item(9,4).Value = "11.97"
$newWs.Cells.Item(9,5).NumberFormat = "@"
$newWs.Cells.Item(9,5).Value = "90.70"
$newWs.Cells.Item(9,6).NumberFormat = "@"
$newWs.Cells.Item(9,6).Value = "4.51"
$newWs.Cells.Item(9,7).NumberFormat = "@"
$newWs.Cells.Item(9,7).Value = "0.5398"
$newWs.Cells.Item(9,8).Value = 7
$newWs.Cells.Item(10,1).Value = 8
$newWs.Cells.Item(10,2).NumberFormat = "@"
$newWs.Cells.Item(10,2).Value = "014967"
$newWs.Cells.Item(10,3).NumberFormat = "@"
$newWs.Cells.Item(10,3).Value = "建信潜力新蓝筹股票C"
$newWs.Cells.Item(10,4).NumberFormat = "@"
$newWs.Cells.Item(10,4).Value = "10.68"
$newWs.Cells.Item(10,5).NumberFormat = "@"
$newWs.Cells.Item(10,5).Value = "90.58"
$newWs.Cells.Item(10,6).NumberFormat = "@"
$newWs.Cells.Item(10,6).Value = "3.75"
$newWs.Cells.Item(10,7).NumberFormat = "@"
$newWs.Cells.Item(10,7).Value = "0.4005"
$newWs.Cells.Item(10,8).Value = 6
$newWs.Cells.Item(11,1).Value = 9
$newWs.Cells.Item(11,2).NumberFormat = "@"
$newWs.Cells.Item(11,2).Value = "013919"
$newWs.Cells.Item(11,3).NumberFormat = "@"
$newWs.Cells.Item(11,3).Value = "建信中小盘先锋股票C"
$newWs.Cells.Item(11,4).NumberFormat = "@"
$newWs.Cells.Item(11,4).Value = "10.14"
$newWs.Cells.Item(11,5).NumberFormat = "@"
$newWs.Cells.Item(11,5).Value = "91.13"
$newWs.Cells.Item(11,6).NumberFormat = "@"
$newWs.Cells.Item(11,6).Value = "3.70"
$newWs.Cells.Item(11,7).NumberFormat = "@"
$newWs.Cells.Item(11,7).Value = "0.3752"
$newWs.Cells.Item(11,8).Value = 6
$newWs.Cells.Item(12,1).Value = 10
$newWs.Cells.Item(12,2).NumberFormat = "@"
$newWs.Cells.Item(12,2).Value = "010452"
$newWs.Cells.Item(12,3).NumberFormat = "@"
$newWs.Cells.Item(12,3).Value = "广发瑞福精选混合A"
$newWs.Cells.Item(12,4).NumberFormat = "@"
$newWs.Cells.Item(12,4).Value = "10.71"
$newWs.Cells.Item(12,5).NumberFormat = "@"
$newWs.Cells.Item(12,5).Value = "83.18"
$newWs.Cells.Item(12,6).NumberFormat = "@"
$newWs.Cells.Item(12,6).Value = "3.00"
$newWs.Cells.Item(12,7).NumberFormat = "@"
$newWs.Cells.Item(12,7).Value = "0.3213"
$newWs.Cells.Item(12,8).Value = 6
$newWs.Cells.Item(13,1).Value = 11
$newWs.Cells.Item(13,2).NumberFormat = "@"
$newWs.Cells.Item(13,2).Value = "013887"
$newWs.Cells.Item(13,3).NumberFormat = "@"
$newWs.Cells.Item(13,3).Value = "华商新能源汽车混合C"
$newWs.Cells.Item(13,4).NumberFormat = "@"
$newWs.Cells.Item(13,4).Value = "3.70"
$newWs.Cells.Item(13,5).NumberFormat = "@"
$newWs.Cells.Item(13,5).Value = "89.00"
$newWs.Cells.Item(13,6).NumberFormat = "@"
$newWs.Cells.Item(13,6).Value = "6.93"
$newWs.Cells.Item(13,7).NumberFormat = "@"
$newWs.Cells.Item(13,7).Value = "0.2564"
$newWs.Cells.Item(13,8).Value = 4
$newWs.Cells.Item(14,1).Value = 12
$newWs.Cells.Item(14,2).NumberFormat = "@"
$newWs.Cells.Item(14,2).Value = "014350"
$newWs.Cells.Item(14,3).NumberFormat = "@"
$newWs.Cells.Item(14,3).Value = "华商卓越成长一年持有混合A"
$newWs.Cells.Item(14,4).NumberFormat = "@"
$newWs.Cells.Item(14,4).Value = "3.05"
$newWs.Cells.Item(14,5).NumberFormat = "@"
$newWs.Cells.Item(14,5).Value = "93.10"
$newWs.Cells.Item(14,6).NumberFormat = "@"
$newWs.Cells.Item(14,6).Value = "5.27"
$newWs.Cells.Item(14,7).NumberFormat = "@"
$newWs.Cells.Item(14,7).Value = "0.1607"
$newWs.Cells.Item(14,8).Value = 5
$newWs.Cells.Item(15,1).Value = 13
$newWs.Cells.Item(15,2).NumberFormat = "@"
$newWs.Cells.Item(15,2).Value = "010453"
$newWs.Cells.Item(15,3).NumberFormat = "@"
$newWs.Cells.Item(15,3).Value = "广发瑞福精选混合C"
$newWs.Cells.Item(15,4).NumberFormat = "@"
$newWs.Cells.Item(15,4).Value = "0.52"
$newWs.Cells.Item(15,5).NumberFormat = "@"
$newWs.Cells.Item(15,5).Value = "83.18"
$newWs.Cells.Item(15,6).NumberFormat = "@"
$newWs.Cells.Item(15,6).Value = "3.00"
$newWs.Cells.Item(15,7).NumberFormat = "@"
$newWs.Cells.Item(15,7).Value = "0.0156"
$newWs.Cells.Item(15,8).Value = 6
$newWs.Cells.Item(16,1).Value = 14
$newWs.Cells.Item(16,2).NumberFormat = "@"
$newWs.Cells.Item(16,2).Value = "014351"
$newWs.Cells.Item(16,3).NumberFormat = "@"
$newWs.Cells.Item(16,3).Value = "华商卓越成长一年持有混合C"
$newWs.Cells.Item(16,4).NumberFormat = "@"
$newWs.Cells.Item(16,4).Value = "0.10"
$newWs.Cells.Item(16,5).NumberFormat = "@"
$newWs.Cells.Item(16,5).Value = "93.10"
$newWs.Cells.Item(16,6).NumberFormat = "@"
$newWs.Cells.Item(16,6).Value = "5.27"
$newWs.Cells.Item(16,7).NumberFormat = "@"
$newWs.Cells.Item(16,7).Value = "0.0053"
$newWs.Cells.Item(16,8).Value = 5
$newWs.Cells.Item(17,1).Value = 15
$newWs.Cells.Item(17,2).NumberFormat = "@"
$newWs.Cells.Item(17,2).Value = "015436"
$newWs.Cells.Item(17,3).NumberFormat = "@"
$newWs.Cells.Item(17,3).Value = "建信优化配置混合C"
$newWs.Cells.Item(17,4).NumberFormat = "@"
$newWs.Cells.Item(17,4).Value = "0.11"
$newWs.Cells.Item(17,5).NumberFormat = "@"
$newWs.Cells.Item(17,5).Value = "88.67"
$newWs.Cells.Item(17,6).NumberFormat = "@"
$newWs.Cells.Item(17,6).Value = "3.60"
$newWs.Cells.Item(17,7).NumberFormat = "@"
$newWs.Cells.Item(17,7).Value = "0.0040"
$newWs.Cells.Item(17,8).Value = 6

# ---------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: insert a new row 2 for the
#    2022-Q4 totals, pushing the existing 2022-Q3 / 2022-Q2 rows down,
#    and renumber the index column (A) sequentially.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Restore plain (unstyled) formatting on the new row from neighbouring
# cells, then fill in the 2022-Q4 totals.
$summary.Cells.Item(3,1).Copy($summary.Cells.Item(2,1))
$summary.Cells.Item(4,2).Copy($summary.Cells.Item(2,2))
$summary.Cells.Item(4,3).Copy($summary.Cells.Item(2,3))
$summary.Cells.Item(4,4).Copy($summary.Cells.Item(2,4))

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).NumberFormat = "@"
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 16
$summary.Cells.Item(2,4).Value = 8.23

# Renumber the existing rows (previously 0/1, now 1/2).
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2

# ---------------------------------------------------------------------
# 4) Restore the original active sheet (2022-Q2 was the selected tab
#    before the edit).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
